$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 19608462
$ws.Range("I33").Value = 33334058
$ws.Range("K33").Value = 33334058
$ws.Range("M33").Value = -33333829
$ws.Range("H40").Value = 1868.6086
$ws.Range("J40").Value = 1438.8334
$ws.Range("L40").Value = 1438.8334
$ws.Range("N40").Value = -1788.8334
$ws.Range("H62").Value = 2836.5715
$ws.Range("I62").Value = 2550
$ws.Range("J62").Value = 2951.2
$ws.Range("K62").Value = 2550
$ws.Range("L62").Value = 2951.2
$ws.Range("M62").Value = -1926
$ws.Range("N62").Value = -4199.2
$ws.Range("H65").Value = 2836.5715
$ws.Range("I65").Value = 2550
$ws.Range("J65").Value = 2951.2
$ws.Range("K65").Value = 12750
$ws.Range("L65").Value = 14756
$ws.Range("M65").Value = -9630
$ws.Range("N65").Value = -20996
$ws.Range("H116").Value = 5198
$ws.Range("I116").Value = 1997.5
$ws.Range("J116").Value = 7331.6665
$ws.Range("K116").Value = 1997.5
$ws.Range("L116").Value = 7331.6665
$ws.Range("M116").Value = 1444.5
$ws.Range("N116").Value = -14215.6665
$ws.Range("H135").Value = 3836.738
$ws.Range("I135").Value = 585.96875
$ws.Range("J135").Value = 14239.2
$ws.Range("K135").Value = 5273.71875
$ws.Range("L135").Value = 128152.8
$ws.Range("M135").Value = -2738.71875
$ws.Range("N135").Value = -133222.8
$ws.Range("H137").Value = 55557670
$ws.Range("I137").Value = 1666.0769
$ws.Range("J137").Value = 200003280
$ws.Range("K137").Value = 4998.2307
$ws.Range("L137").Value = 600009840
$ws.Range("M137").Value = -2448.2307
$ws.Range("N137").Value = -600014940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H45").Value = 15152397
$ws.Range("I45").Value = 22222962
$ws.Range("J45").Value = 1185.7142
$ws.Range("K45").Value = 22222962
$ws.Range("L45").Value = 1185.7142
$ws.Range("M45").Value = -22222585
$ws.Range("N45").Value = -1939.7142
$ws.Range("H61").Value = 2113.44
$ws.Range("I61").Value = 1997.3
$ws.Range("K61").Value = 1997.3
$ws.Range("M61").Value = -1785.3
$ws.Range("H74").Value = 1024.3903
$ws.Range("I74").Value = 1017.4167
$ws.Range("J74").Value = 1074.6
$ws.Range("K74").Value = 1017.4167
$ws.Range("L74").Value = 1074.6
$ws.Range("M74").Value = -143.4167
$ws.Range("N74").Value = -2822.6
$ws.Range("H77").Value = 1024.3903
$ws.Range("I77").Value = 1017.4167
$ws.Range("J77").Value = 1074.6
$ws.Range("K77").Value = 5087.0835
$ws.Range("L77").Value = 5373
$ws.Range("M77").Value = -719.0834999999997
$ws.Range("N77").Value = -14109
$ws.Range("H132").Value = 6222.38
$ws.Range("I132").Value = 6369.1904
$ws.Range("J132").Value = 5451.625
$ws.Range("K132").Value = 19107.5712
$ws.Range("L132").Value = 16354.875
$ws.Range("M132").Value = -16577.5712
$ws.Range("N132").Value = -21414.875
$ws.Range("H136").Value = 2113.44
$ws.Range("I136").Value = 1997.3
$ws.Range("K136").Value = 5991.9
$ws.Range("M136").Value = -3441.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 7462.222
$ws.Range("J81").Value = 7462.222
$ws.Range("L81").Value = 7462.222
$ws.Range("N81").Value = -9584.222
$ws.Range("H84").Value = 7462.222
$ws.Range("J84").Value = 7462.222
$ws.Range("L84").Value = 22386.666
$ws.Range("N84").Value = -32994.666
$ws.Range("H99").Value = 749.8570999999999
$ws.Range("I99").Value = 741.6667
$ws.Range("J99").Value = 799
$ws.Range("K99").Value = 741.6667
$ws.Range("L99").Value = 799
$ws.Range("M99").Value = 756.3333
$ws.Range("N99").Value = -3795
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H58").Value = 1394.5676
$ws.Range("I58").Value = 1434.4572
$ws.Range("J58").Value = 696.5
$ws.Range("K58").Value = 1434.4572
$ws.Range("L58").Value = 696.5
$ws.Range("M58").Value = -1231.4572
$ws.Range("N58").Value = -1102.5
$ws.Range("H99").Value = 2499.9412
$ws.Range("I99").Value = 3174.75
$ws.Range("K99").Value = 3174.75
$ws.Range("M99").Value = -1676.75
$ws.Range("H107").Value = 704.1852
$ws.Range("I107").Value = 679.82855
$ws.Range("K107").Value = 679.82855
$ws.Range("M107").Value = 1240.17145
$ws.Range("H126").Value = 2499.9412
$ws.Range("I126").Value = 3174.75
$ws.Range("K126").Value = 9524.25
$ws.Range("M126").Value = -7054.25
$ws.Range("H132").Value = 4033602.2
$ws.Range("I132").Value = 1053.9259
$ws.Range("J132").Value = 31253304
$ws.Range("K132").Value = 3161.7777
$ws.Range("L132").Value = 93759912
$ws.Range("M132").Value = -631.7776999999996
$ws.Range("N132").Value = -93764972
$ws.Range("H134").Value = 2578.32
$ws.Range("I134").Value = 2578.32
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7734.960000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5199.960000000001
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1394.5676
$ws.Range("I136").Value = 1434.4572
$ws.Range("J136").Value = 696.5
$ws.Range("K136").Value = 4303.3716
$ws.Range("L136").Value = 2089.5
$ws.Range("M136").Value = -1753.3716
$ws.Range("N136").Value = -7189.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 221.4
$ws.Range("I6").Value = 26.75
$ws.Range("K6").Value = 80.25
$ws.Range("M6").Value = 32.75
$ws.Range("H12").Value = 122.80769
$ws.Range("I12").Value = 95.333336
$ws.Range("K12").Value = 286.000008
$ws.Range("M12").Value = -113.000008
$ws.Range("H117").Value = 10785.692
$ws.Range("I117").Value = 887.6
$ws.Range("K117").Value = 2662.8
$ws.Range("M117").Value = 779.1999999999998
$ws.Range("H131").Value = 2099057.2
$ws.Range("I131").Value = 11495.444
$ws.Range("J131").Value = 2526058.5
$ws.Range("K131").Value = 34486.33199999999
$ws.Range("L131").Value = 7578175.5
$ws.Range("M131").Value = -29446.33199999999
$ws.Range("N131").Value = -7588255.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 14967.25
$ws.Range("I3").Value = 33583.332
$ws.Range("J3").Value = 3797.6
$ws.Range("K3").Value = 33583.332
$ws.Range("L3").Value = 3797.6
$ws.Range("M3").Value = -33467.332
$ws.Range("N3").Value = -4029.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10962.096
$ws.Range("I132").Value = 16400.666
$ws.Range("J132").Value = 3710.6667
$ws.Range("K132").Value = 49201.99800000001
$ws.Range("L132").Value = 11132.0001
$ws.Range("M132").Value = -46671.99800000001
$ws.Range("N132").Value = -16192.0001
$ws.Range("H136").Value = 5166.1177
$ws.Range("I136").Value = 5401.6553
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 16204.9659
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -13654.9659
$ws.Range("N136").Value = -16500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 41667380
$ws.Range("I96").Value = 125000390
$ws.Range("J96").Value = 874.75
$ws.Range("K96").Value = 125000390
$ws.Range("L96").Value = 874.75
$ws.Range("M96").Value = -124999017
$ws.Range("N96").Value = -3620.75
$ws.Range("H107").Value = 12387741
$ws.Range("I107").Value = 5000348.5
$ws.Range("J107").Value = 27778144
$ws.Range("K107").Value = 15001045.5
$ws.Range("L107").Value = 83334432
$ws.Range("M107").Value = -14999125.5
$ws.Range("N107").Value = -83338272
$ws.Range("H132").Value = 1974.7407
$ws.Range("I132").Value = 1732.8
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5198.4
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2668.4
$ws.Range("N132").Value = -20057
